# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style from H1, and fill in the data rows 2-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(7,9,6,8,9,8,9,7,7,7,9,8,8,8,6,4,7,7,7,7,8,7,8,6,6,10,5,5,7,8,6,8,6,9,8,7,7,6,9,7,9,7,6,6,10,7,6,6,7,5,7,7,7,7,6,7,3,3)
$jValues = @(8,9,7,9,9,9,9,7,7,7,9,8,8,8,7,7,7,7,7,7,8,9,8,6,6,10,6,5,7,8,6,8,6,9,8,7,7,7,9,7,9,7,6,6,10,7,7,6,7,6,7,7,7,8,6,7,3,3)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
